# Insert a new weekly price record at row 35 for "Vega Monumental Concepción - Berenjena".
# All subsequent rows (old 35..61) shift down by one (to 36..62), which Excel's
# native row-insert semantics already take care of (including carrying the
# date-column number format down from the row above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 35, pushing existing rows 35-61 to 36-62.
$ws.Rows(35).Insert()

# Populate the new row 35 with the new weekly entry.
$ws.Range("A35").Value2 = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value2 = 44574
$ws.Range("E35").Value2 = 8
$ws.Range("F35").Value2 = 100112001
$ws.Range("G35").Value = "Berenjena"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value2 = 100
$ws.Range("K35").Value2 = 10000
$ws.Range("L35").Value2 = 11000
$ws.Range("M35").Value2 = 10500
$ws.Range("N35").Value = "$/caja 60 unidades"
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value2 = 175
$ws.Range("Q35").Value2 = 60
$ws.Range("R35").Value = "Hortaliza"
